# Applies: add UNIQUE_ID column to TEAM_PLAYER_MAPPINGS, shift old columns right,
# update active tab/selection on several sheets (scoring calculation logic).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) TEAM_PLAYER_MAPPINGS (sheet2): insert a UNIQUE_ID column at C, shifting
#    the existing PLAYER_ROLE (old C) -> D and TEAM_INITIALS (old D) -> E.
# ---------------------------------------------------------------------------
$wsMap = $wb.Worksheets.Item("TEAM_PLAYER_MAPPINGS")

# Shift existing columns D->E and C->D first (read with Value2 to avoid
# picking up values we are about to overwrite).
for ($r = 1; $r -le 20; $r++) {
    $oldD = $wsMap.Cells.Item($r, 4).Value2
    $wsMap.Cells.Item($r, 5).Value = $oldD
}
for ($r = 1; $r -le 20; $r++) {
    $oldC = $wsMap.Cells.Item($r, 3).Value2
    $wsMap.Cells.Item($r, 4).Value = $oldC
}

# Now populate the new UNIQUE_ID column C, in row order (header then 19
# players) so new shared-string entries are created in the expected order.
$wsMap.Cells.Item(1, 3).Value = "UNIQUE_ID"
$wsMap.Cells.Item(2, 3).Value = "RS"
$wsMap.Cells.Item(3, 3).Value = "HS"
$wsMap.Cells.Item(4, 3).Value = "LM"
$wsMap.Cells.Item(5, 3).Value = "KP"
$wsMap.Cells.Item(6, 3).Value = "MSD"
$wsMap.Cells.Item(7, 3).Value = "RJ"
$wsMap.Cells.Item(8, 3).Value = "DB"
$wsMap.Cells.Item(9, 3).Value = "VK"
$wsMap.Cells.Item(10, 3).Value = "ABdeV"
$wsMap.Cells.Item(11, 3).Value = "SW"
$wsMap.Cells.Item(12, 3).Value = "KLR"
$wsMap.Cells.Item(13, 3).Value = "GG"
$wsMap.Cells.Item(14, 3).Value = "YP"
$wsMap.Cells.Item(15, 3).Value = "RU"
$wsMap.Cells.Item(16, 3).Value = "SN"
$wsMap.Cells.Item(17, 3).Value = "DW"
$wsMap.Cells.Item(18, 3).Value = "SD"
$wsMap.Cells.Item(19, 3).Value = "BK"
$wsMap.Cells.Item(20, 3).Value = "AN"

[void]$wsMap.Range("C8").Select()

# ---------------------------------------------------------------------------
# 2) LEAGUE_RULES (sheet3): move selection to B4.
# ---------------------------------------------------------------------------
$wsRules = $wb.Worksheets.Item("LEAGUE_RULES")
[void]$wsRules.Range("B4").Select()

# ---------------------------------------------------------------------------
# 3) GAMES (sheet4): move selection to D3.
# ---------------------------------------------------------------------------
$wsGames = $wb.Worksheets.Item("GAMES")
[void]$wsGames.Range("D3").Select()

# ---------------------------------------------------------------------------
# 4) Make TEAM_PLAYER_MAPPINGS the active sheet/tab (activeTab index 1, and
#    tabSelected on sheet2 while GAMES loses tabSelected).
# ---------------------------------------------------------------------------
[void]$wsMap.Activate()
[void]$wsMap.Range("C8").Select()
